# Yash_DeploymentSheet.xlsx update
# - Adds two new deployment-parameter rows (RPAChallenge_Path / RPAChallenge_URL)
# - Bumps the package version from 1.0.1 to 1.0.2
# - Widens column C to fit the long path value
# - Leaves the final selection on C23 (matches the saved workbook state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values first, in the same order the original author typed them,
#     so the shared-string table comes out in the same sequence. ---
$ws.Range("A9").Value  = "RPAChallenge_Path"
$ws.Range("B9").Value  = "Text"
$ws.Range("A10").Value = "RPAChallenge_URL"
$ws.Range("C9").Value  = "C:\Users\User\OneDrive\Documents\RPA30Days_Program\RPA-Developer-in-30-Days\SourceCode\2023\August2023\Yash\RE_RPAChallenge\Data\Input\challenge.xlsx"
$ws.Range("B10").Value = "Text"
$ws.Range("C10").Value = "https://rpachallenge.com/"
$ws.Range("C17").Value = "1.0.2"

# --- Formatting: A9/A10 use a bold, dark-grey Calibri label style ---
$ws.Range("A9").Font.Bold  = $true
$ws.Range("A9").Font.Color = 5590598   # RGB(0x46,0x4E,0x55) == FF464E55 in BGR order

$ws.Range("A10").Font.Bold  = $true
$ws.Range("A10").Font.Color = 5590598

# --- Formatting: B10 uses a lighter grey Segoe UI style ---
$ws.Range("B10").Font.Name  = "Segoe UI"
$ws.Range("B10").Font.Color = 7300699  # RGB(0x5B,0x66,0x6F) == FF5B666F in BGR order

# --- Row 10 is visually taller, like row 23 already is ---
$ws.Rows.Item(10).RowHeight = 16.8

# --- Column C must be widened to fit the long file-path text ---
$ws.Columns.Item(3).ColumnWidth = 140.43

# --- Final selection lands on C23 ---
$ws.Range("C23").Select() | Out-Null

Write-Output "Deployment sheet updated: RPAChallenge_Path/RPAChallenge_URL rows added, version bumped to 1.0.2"
